$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append rows 10-19 (Id 9-18) to the product list: a sequential Id in
# column A, and empty (but present) text cells in B-E, matching the
# "vaciado de lista" (empty cart-window placeholder rows) added by the
# commit.
for ($i = 10; $i -le 19; $i++) {
    $id = $i - 1
    $ws.Cells.Item($i, 1).Value = $id

    foreach ($col in 2..5) {
        $cell = $ws.Cells.Item($i, $col)
        $cell.Formula = "'"       # force an explicit empty-text cell
        $cell.ClearFormats()      # drop the quote-prefix style bit again
    }
}
